# Commit: set/get auction on group
# This script applies the edits to "CricDream issue list.xlsx" (Sheet1)
# that correspond to the target diff:
#  - Issues #5, #6, #13, #17 (rows 7, 8, 15, 19) get marked/highlighted
#    (yellow fill) and their Status is changed from "Pending" to "Resolved".
#  - Issue #7 (row 9, about group-only purchase) gets a Resolved Description
#    of "Not to be done. Ignore".
#  - The active selection moves to D19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$yellow = 65535  # RGB(255,255,0) as OLE color (0x00FFFF00 -> BGR 65535)

# Rows whose Issue Description (column B) gets highlighted and whose
# Status (column D) flips from Pending to Resolved.
$rows = @(7, 8, 15, 19)
foreach ($r in $rows) {
    $ws.Range("B$r").Interior.Color = $yellow
    $ws.Range("D$r").Value = "Resolved"
}

# Row 9 only gets its Status resolved plus a Resolved Description note.
$ws.Range("D9").Value = "Resolved"
$ws.Range("F9").Value = "Not to be done. Ignore"

# Restore selection to match the saved view state of the edited workbook.
$ws.Range("D19").Select()
